$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 2953.32
$ws.Range("I113").Value = 2633.4375
$ws.Range("J113").Value = 3522
$ws.Range("K113").Value = 2633.4375
$ws.Range("L113").Value = 3522
$ws.Range("M113").Value = 620.5625
$ws.Range("N113").Value = -10030

$ws.Range("H132").Value = 4720012
$ws.Range("I132").Value = 2630.7234
$ws.Range("K132").Value = 7892.1702
$ws.Range("M132").Value = -5362.1702

$ws.Range("H138").Value = 4699602
$ws.Range("I138").Value = 14495793
$ws.Range("J138").Value = 5593.375
$ws.Range("K138").Value = 43487379
$ws.Range("L138").Value = 16780.125
$ws.Range("M138").Value = -43482239
$ws.Range("N138").Value = -27060.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()

$ws.Range("H32").Value = 15381.241
$ws.Range("I32").Value = 15146.787
$ws.Range("J32").Value = 21633.334
$ws.Range("K32").Value = 15146.787
$ws.Range("L32").Value = 21633.334
$ws.Range("M32").Value = -14859.787
$ws.Range("N32").Value = -22207.334

$ws.Range("H45").Value = 1378845.9
$ws.Range("I45").Value = 2274170.2
$ws.Range("J45").Value = 1423.4615
$ws.Range("K45").Value = 2274170.2
$ws.Range("L45").Value = 1423.4615
$ws.Range("M45").Value = -2273793.2
$ws.Range("N45").Value = -2177.4615

$ws.Range("H61").Value = 1910.9474
$ws.Range("I61").Value = 1756.48
$ws.Range("J61").Value = 3014.2856
$ws.Range("K61").Value = 1756.48
$ws.Range("L61").Value = 3014.2856
$ws.Range("M61").Value = -1544.48
$ws.Range("N61").Value = -3438.2856

$ws.Range("H74").Value = 4534.3667
$ws.Range("I74").Value = 1048.2
$ws.Range("J74").Value = 21965.2
$ws.Range("K74").Value = 1048.2
$ws.Range("L74").Value = 21965.2
$ws.Range("M74").Value = -174.2
$ws.Range("N74").Value = -23713.2

$ws.Range("H77").Value = 4534.3667
$ws.Range("I77").Value = 1048.2
$ws.Range("J77").Value = 21965.2
$ws.Range("K77").Value = 5241
$ws.Range("L77").Value = 109826
$ws.Range("M77").Value = -873
$ws.Range("N77").Value = -118562

$ws.Range("H136").Value = 1910.9474
$ws.Range("I136").Value = 1756.48
$ws.Range("J136").Value = 3014.2856
$ws.Range("K136").Value = 5269.440000000001
$ws.Range("L136").Value = 9042.856800000001
$ws.Range("M136").Value = -2719.440000000001
$ws.Range("N136").Value = -14142.8568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 1214.5
$ws.Range("I22").Value = 1119.3334
$ws.Range("J22").Value = 1500
$ws.Range("K22").Value = 1119.3334
$ws.Range("L22").Value = 1500
$ws.Range("M22").Value = -946.3334
$ws.Range("N22").Value = -1846

$ws.Range("H64").Value = 634.55554
$ws.Range("J64").Value = 833
$ws.Range("L64").Value = 833
$ws.Range("N64").Value = -1283

$ws.Range("H67").Value = 634.55554
$ws.Range("J67").Value = 833
$ws.Range("L67").Value = 833
$ws.Range("N67").Value = -2393

$ws.Range("H99").Value = 607.5
$ws.Range("I99").Value = 607.5
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 607.5
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 890.5
$ws.Range("N99").ClearContents()

$ws.Range("H118").Value = 7875.5557
$ws.Range("J118").Value = 7875.5557
$ws.Range("L118").Value = 7875.5557
$ws.Range("N118").Value = -11189.5557

$ws.Range("H132").Value = 40822.418
$ws.Range("J132").Value = 40822.418
$ws.Range("L132").Value = 40822.418
$ws.Range("N132").Value = -50942.418

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4968.7383
$ws.Range("I31").Value = 1232.5667
$ws.Range("J31").Value = 8171.1714
$ws.Range("K31").Value = 1232.5667
$ws.Range("L31").Value = 8171.1714
$ws.Range("M31").Value = -937.5667000000001
$ws.Range("N31").Value = -8761.171399999999

$ws.Range("H34").Value = 4968.7383
$ws.Range("I34").Value = 1232.5667
$ws.Range("J34").Value = 8171.1714
$ws.Range("K34").Value = 1232.5667
$ws.Range("L34").Value = 8171.1714
$ws.Range("M34").Value = -1030.5667
$ws.Range("N34").Value = -8575.171399999999

$ws.Range("H58").Value = 3102
$ws.Range("I58").Value = 1641.2307
$ws.Range("J58").Value = 6900
$ws.Range("K58").Value = 1641.2307
$ws.Range("L58").Value = 6900
$ws.Range("M58").Value = -1438.2307
$ws.Range("N58").Value = -7306

$ws.Range("H105").Value = 2511.7778
$ws.Range("I105").Value = 2501.1428
$ws.Range("J105").Value = 2549
$ws.Range("K105").Value = 2501.1428
$ws.Range("L105").Value = 2549
$ws.Range("M105").Value = -754.1428000000001
$ws.Range("N105").Value = -6043

$ws.Range("H107").Value = 407.0909
$ws.Range("I107").Value = 365.3125
$ws.Range("J107").Value = 518.5
$ws.Range("K107").Value = 365.3125
$ws.Range("L107").Value = 518.5
$ws.Range("M107").Value = 1554.6875
$ws.Range("N107").Value = -4358.5

$ws.Range("H136").Value = 3102
$ws.Range("I136").Value = 1641.2307
$ws.Range("J136").Value = 6900
$ws.Range("K136").Value = 4923.6921
$ws.Range("L136").Value = 20700
$ws.Range("M136").Value = -2373.6921
$ws.Range("N136").Value = -25800

$ws.Range("H137").Value = 64427.5
$ws.Range("J137").Value = 64103.332
$ws.Range("L137").Value = 64103.332
$ws.Range("N137").Value = -74303.33199999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2702.7144
$ws.Range("I131").Value = 4099.909
$ws.Range("J131").Value = 2062.3333
$ws.Range("K131").Value = 12299.727
$ws.Range("L131").Value = 6186.999899999999
$ws.Range("M131").Value = -7259.726999999999
$ws.Range("N131").Value = -16266.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 7786
$ws.Range("I102").Value = 11478
$ws.Range("K102").Value = 11478
$ws.Range("M102").Value = -9856

$ws.Range("H132").Value = 7397.591
$ws.Range("I132").Value = 10152.5
$ws.Range("J132").Value = 2576.5
$ws.Range("K132").Value = 30457.5
$ws.Range("L132").Value = 7729.5
$ws.Range("M132").Value = -27927.5
$ws.Range("N132").Value = -12789.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7955.32
$ws.Range("I7").Value = 5875.769
$ws.Range("J7").Value = 10208.167
$ws.Range("K7").Value = 5875.769
$ws.Range("L7").Value = 10208.167
$ws.Range("M7").Value = -5763.769
$ws.Range("N7").Value = -10432.167

$ws.Range("H46").Value = 1680
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 1680
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 1680
$ws.Range("N46").Value = -2056
$ws.Range("M46").ClearContents()

$ws.Range("H126").Value = 7955.32
$ws.Range("I126").Value = 5875.769
$ws.Range("J126").Value = 10208.167
$ws.Range("K126").Value = 17627.307
$ws.Range("L126").Value = 30624.501
$ws.Range("M126").Value = -15157.307
$ws.Range("N126").Value = -35564.501

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 15355.772
$ws.Range("J64").Value = 15355.772
$ws.Range("L64").Value = 15355.772
$ws.Range("N64").Value = -15851.772

$ws.Range("H67").Value = 15355.772
$ws.Range("J67").Value = 15355.772
$ws.Range("L67").Value = 15355.772
$ws.Range("N67").Value = -17071.772

$ws.Range("H136").Value = 1268.9615
$ws.Range("I136").Value = 923.27905
$ws.Range("J136").Value = 2920.5557
$ws.Range("K136").Value = 2769.83715
$ws.Range("L136").Value = 8761.667099999999
$ws.Range("M136").Value = -219.8371499999998
$ws.Range("N136").Value = -13861.6671
